$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark (bookmarkStart/bookmarkEnd pair) that
#    followed the "(only ratings of 4 and 5)" run.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Update the "Scrutability" bullet text and split the run so the word
#    "perceive" (replacing "tell") lives in its own run, matching the
#    structure produced by Word when text is retyped mid-sentence.
$d.Content.Find.Execute("Scrutability: allow users to tell if system is mistaken.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Scrutability: allow users to perceive if system is mistaken.", 2)

$r = $d.Content
$r.Find.Execute("perceive", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Bold = $true
$r.Bold = $false
